$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# rewardMoney_sub (column J) was blank for these rows; fill with "0z"
$jRows = @(2,3,5,6,7,8,11,20,22,23,27,28,41,49,50,54,55,56,69,76,82,126)
foreach ($r in $jRows) {
    $ws.Cells.Item($r, 10).Value = "0z"
}

# down_payment (column H) was blank for these rows; fill with "0z"
$hRows = @(164,169,191,194,202,216)
foreach ($r in $hRows) {
    $ws.Cells.Item($r, 8).Value = "0z"
}
